$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 3 entirely; this shifts all rows below it up by one,
# matching the "import maks 10000 row berikutnya" update where the
# first listed row (size 37 for DL00023-1LUBCK) is dropped and the
# remaining rows shift into its place.
$ws.Rows.Item(3).Delete()

# Update the active cell selection to A3, as in the saved workbook.
$ws.Range("A3").Select()
